$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 233
$ws.Range("B3").Value = 162

# Add new rows 4 and 5, copying the style of A3 (style index 1) into A4/A5
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 127
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 108
